# paises.xlsx — "Update countries & provincias Spain"
#
# The underlying country table (sheet "Pais") is kept sorted by "Casos
# totales" (column B) descending. New data came in between the 20:52 and
# 21:22 snapshots, which (a) updated the stats of several countries and
# (b) re-sorted a handful of rows, so some row positions now show a
# different country than before even though their numeric-column values
# simply shifted to the next rank. We therefore:
#   1. bump the "last updated" timestamp in A1,
#   2. re-label the A-column country name for every row whose rank
#      changed, and
#   3. write the new B:H statistics for every affected row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Timestamp -----------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 21:22"

# 2) Country re-labels caused by re-sorting -------------------------------
$countryRelabels = @(
    @("A33", "Ecuador"),
    @("A34", "Peru"),
    @("A35", "Japon"),
    @("A36", "Pakistan"),
    @("A37", "Malasia"),
    @("A38", "Filipinas"),
    @("A39", "Arabia Saudita"),
    @("A40", "Luxemburgo"),
    @("A41", "Indonesia"),
    @("A111", "Georgia"),
    @("A112", "Bolivia"),
    @("A196", "Santa Sede"),
    @("A197", "Sierra Leona"),
    @("A199", "San Bartolome"),
    @("A200", "Nicaragua"),
    @("A208", "Anguila"),
    @("A209", "Burundi"),
    @("A210", "Islas Virgenes Britanicas"),
    @("A212", "Bonaire, San Eustaquio y Saba"),
    @("A213", "Papua Nueva Guinea")
)
foreach ($pair in $countryRelabels) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# 3) Updated statistics (Casos totales, Nuevos casos, Casos activos,
#    Recuperados, Casos criticos, Muertes hoy, Muertes) ------------------
$statUpdates = @(
    @("B4", 420705), @("C4", 20370), @("D4", 22187), @("E4", 384149), @("G4", 1528), @("H4", 14369),
    @("B16", 19183), @("C16", 1286), @("D16", 4533), @("E16", 14223),
    @("B19", 12937), @("C19", 298), @("E19", 8152),
    @("E21", 8530), @("G21", 8), @("H21", 73),
    @("B27", 5916), @("C27", 565), @("E27", 5232),
    @("B33", 4450), @("C33", 455), @("D33", 140), @("E33", 4068), @("F33", 146), @("G33", 22), @("H33", 242),
    @("B34", 4342), @("C34", 1388), @("D34", 1301), @("E34", 2934), @("F34", 109), @("G34", 0), @("H34", 107),
    @("B35", 4257), @("C35", 0), @("D35", 622), @("E35", 3542), @("F35", 80), @("G35", 0), @("H35", 93),
    @("B36", 4196), @("C36", 161), @("D36", 467), @("E36", 3669), @("F36", 25), @("G36", 3), @("H36", 60),
    @("B37", 4119), @("C37", 156), @("D37", 1487), @("E37", 2567), @("F37", 76), @("G37", 2), @("H37", 65),
    @("B38", 3870), @("C38", 106), @("D38", 96), @("E38", 3592), @("F38", 1), @("G38", 5), @("H38", 182),
    @("B39", 3122), @("C39", 327), @("D39", 631), @("E39", 2450), @("F39", 41), @("G39", 0), @("H39", 41),
    @("B40", 3034), @("C40", 64), @("D40", 500), @("E40", 2488), @("F40", 34), @("G40", 2), @("H40", 46),
    @("B41", 2956), @("C41", 218), @("D41", 222), @("E41", 2494), @("F41", 0), @("G41", 19), @("H41", 240),
    @("B75", 727), @("C75", 30), @("E75", 666),
    @("B83", 576), @("C83", 28), @("E83", 495),
    @("B87", 502), @("C87", 19), @("D87", 29), @("E87", 471), @("F87", 15),
    @("B111", 211), @("C111", 15), @("D111", 50), @("E111", 158), @("F111", 6), @("G111", 0), @("H111", 3),
    @("B112", 210), @("C112", 16), @("D112", 2), @("E112", 193), @("F112", 3), @("G112", 1), @("H112", 15),
    @("E153", 19), @("G153", 1), @("H153", 6),
    @("B196", 8), @("D196", 2), @("E196", 6),
    @("C197", 1),
    @("D199", 1), @("H199", 0),
    @("D200", 0), @("H200", 1)
)
foreach ($pair in $statUpdates) {
    $ws.Range($pair[0]).Value = $pair[1]
}
